# Actualización automática 2025-07-21 09:20:08
#
# Registers a new PORCELANATO sale of 366.83 for client
# "CAIZA COLLAGUAZO ROCIO PILAR" (RIOS CARRION ANGEL BENIGNO) in julio,
# and propagates the resulting totals / percentages / counters across
# the three sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" --------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Client "CAIZA COLLAGUAZO ROCIO PILAR" (row 5): PORCELANATO sale recorded.
$wsGrupo.Range("M5").Value = 366.83

# Totals row: number of advisors who reached their PORCELANATO goal goes
# from 3 to 4 out of 22.
$wsGrupo.Range("M24").Value = "4 de 22"

# --- Sheet "VENTA MENSUAL" ------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Same sale shows up under julio (column F) for the same client (row 5).
$wsMensual.Range("F5").Value = 366.83

# Totals row for julio increases by the same amount.
$wsMensual.Range("F24").Value = 17480.29

# --- Sheet "CUMPLIMIENTO MENSUAL" ----------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column F (CUMPLIMIENTO) narrows slightly following the layout refresh.
# ColumnWidth is offset from the stored XML width by ~0.8333 characters,
# so 22.1666... renders as the target width of 23 in the saved file.
$wsCumpl.Range("F1").EntireColumn.ColumnWidth = 22.1666666666667

# PORCELANATO group (row 16): VENTA, POR CUMPLIR and CUMPLIMIENTO updated
# to reflect the new sale.
$wsCumpl.Range("D16").Value = 17480.29
$wsCumpl.Range("E16").Value = 21276.25
$wsCumpl.Range("F16").Value = 0.451028136154569

# TOTAL row (row 19): same propagation.
$wsCumpl.Range("D19").Value = 17480.29
$wsCumpl.Range("E19").Value = 40742.71386304603
$wsCumpl.Range("F19").Value = 0.3002299579238076
